$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wydatki budowa")

# Insert two new rows right above the blank spacer row 58 that sits just
# above the "Suma" totals row of the Tabela1 table - this pushes the old
# rows 58-61 (blank spacer rows + totals row + spacer) down by two, like
# using Excel's "Insert Table Rows Above" there.
$ws.Rows("58:59").Insert()

# The new row 58 should look like the invoice line above it (row 57) -
# match that formatting (in particular column A, which carries a subtle
# border tied to the "Etap" grouping) before filling in the values.
$ws.Range("A57").Copy()
$ws.Range("A58").PasteSpecial(-4122)

# Fill in the new invoice line (row 58) - another steel delivery for the
# ceiling slab ("kolejna stal na strop").
$ws.Range("A58").Value = "2. Ściany nadziemia"
$ws.Range("B58").Value = "Materiał"
$ws.Range("C58").Value = "Manex"
$ws.Range("D58").Value = "Stal"
$ws.Range("E58").Value = 227.55
$ws.Range("F58").Value = "5801/T/09/2013"
$ws.Range("G58").Value = "9/11/2013"
$ws.Range("H58").Value = "9/14/2013"
$ws.Range("I58").Value = "9/11/2013"
$ws.Range("J58").Value = "eb wsp"

# Grow the table (Tabela1) so it keeps covering the data through the
# (now shifted) totals row.
$lo = $ws.ListObjects.Item("Tabela1")
$lo.Resize($ws.Range("A1:J62"))

# Keep the selection where the user left it after adding the new row.
$ws.Range("J60").Select()
